# Daily attendance processing - 2025-11-03 20:47:33
#
# Normalizes the "Recorded By" (column G) values on the
# "Session Analysis Results" sheet:
#   - "System, <email>"                      -> "<email>, System"
#   - "System, system, backup@backdoor.com"  -> "system, System, backup@backdoor.com"
# Rows whose "Recorded By" text doesn't match one of the old forms above
# (e.g. "System, backup@backdoor.com", single-name entries, or rows
# already in the target order) are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$renames = @{
    2   = 'system, System, backup@backdoor.com'
    3   = 'dnasr281@gmail.com, System'
    6   = 'dnasr281@gmail.com, System'
    7   = 'admin@admin.com, System'
    10  = 'dnasr281@gmail.com, System'
    12  = 'dnasr281@gmail.com, System'
    13  = 'dnasr281@gmail.com, System'
    14  = 'dnasr281@gmail.com, System'
    15  = 'dnasr281@gmail.com, System'
    18  = 'dnasr281@gmail.com, System'
    19  = 'dnasr281@gmail.com, System'
    20  = 'dnasr281@gmail.com, System'
    21  = 'dnasr281@gmail.com, System'
    22  = 'dnasr281@gmail.com, System'
    24  = 'dnasr281@gmail.com, System'
    26  = 'dnasr281@gmail.com, System'
    28  = 'system, System, backup@backdoor.com'
    29  = 'dnasr281@gmail.com, System'
    32  = 'dnasr281@gmail.com, System'
    33  = 'admin@admin.com, System'
    36  = 'dnasr281@gmail.com, System'
    38  = 'dnasr281@gmail.com, System'
    39  = 'dnasr281@gmail.com, System'
    40  = 'dnasr281@gmail.com, System'
    41  = 'dnasr281@gmail.com, System'
    44  = 'dnasr281@gmail.com, System'
    45  = 'dnasr281@gmail.com, System'
    46  = 'dnasr281@gmail.com, System'
    47  = 'dnasr281@gmail.com, System'
    48  = 'dnasr281@gmail.com, System'
    50  = 'dnasr281@gmail.com, System'
    52  = 'dnasr281@gmail.com, System'
    54  = 'system, System, backup@backdoor.com'
    55  = 'dnasr281@gmail.com, System'
    58  = 'dnasr281@gmail.com, System'
    59  = 'admin@admin.com, System'
    62  = 'dnasr281@gmail.com, System'
    64  = 'dnasr281@gmail.com, System'
    65  = 'dnasr281@gmail.com, System'
    66  = 'dnasr281@gmail.com, System'
    67  = 'dnasr281@gmail.com, System'
    70  = 'dnasr281@gmail.com, System'
    71  = 'dnasr281@gmail.com, System'
    72  = 'dnasr281@gmail.com, System'
    73  = 'dnasr281@gmail.com, System'
    74  = 'dnasr281@gmail.com, System'
    76  = 'dnasr281@gmail.com, System'
    78  = 'dnasr281@gmail.com, System'
    83  = 'dnasr281@gmail.com, System'
    84  = 'dnasr281@gmail.com, System'
    85  = 'dnasr281@gmail.com, System'
    86  = 'dnasr281@gmail.com, System'
    90  = 'dnasr281@gmail.com, System'
    92  = 'dnasr281@gmail.com, System'
    99  = 'dnasr281@gmail.com, System'
    101 = 'dnasr281@gmail.com, System'
    109 = 'dnasr281@gmail.com, System'
    110 = 'dnasr281@gmail.com, System'
    111 = 'dnasr281@gmail.com, System'
    112 = 'dnasr281@gmail.com, System'
    116 = 'dnasr281@gmail.com, System'
    118 = 'dnasr281@gmail.com, System'
    125 = 'dnasr281@gmail.com, System'
    127 = 'dnasr281@gmail.com, System'
    135 = 'dnasr281@gmail.com, System'
    136 = 'dnasr281@gmail.com, System'
    137 = 'dnasr281@gmail.com, System'
    138 = 'dnasr281@gmail.com, System'
    142 = 'dnasr281@gmail.com, System'
    144 = 'dnasr281@gmail.com, System'
    151 = 'dnasr281@gmail.com, System'
    153 = 'dnasr281@gmail.com, System'
}

foreach ($row in $renames.Keys) {
    $ws.Range("G$row").Value = $renames[$row]
}
